$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared text blocks used in the new rows (written in the same order
# they first appear so the shared-strings table gets built in the right
# sequence).
$txt_B209 = "지름 `$\mathrm{AB}`$ 의 중점을 `$\mathrm{M}`$ 이라하고 `$f(\theta)`$를 부채꼴과 삼각형으로 분해해 봅니다."
$txt_B210 = "부채꼴 `$\mathrm{AMQ}`$의 넓이를 `$\theta`$의 식으로 구합니다."
$txt_B211 = "삼각형 `$\mathrm{MBQ}`$ 의 넓이를 `$\theta`$의 식으로 구합니다."
$txt_B212 = "삼각형 `$\mathrm{RAB}`$ 의 넓이를 `$\theta`$의 식으로 구합니다."
$txt_B213 = "`$f(\theta)`$를 구해서 `$\displaystyle\lim _{\theta \rightarrow 0+} \dfrac{f(\theta)}{\theta}`$를 계산합니다."
$txt_B214 = "평행선에 의한 닮음을 이용해서 정삼각형 `$\mathrm{STU}`$ 의 한 변의 길이를 `$\theta`$의 식으로 구합니다."
$txt_B215 = "`$g(\theta)`$를 구해서 `$\displaystyle\lim _{\theta \rightarrow 0+} \dfrac{g(\theta)}{\theta^{2}}`$를 계산합니다."
$txt_B216 = "수렴하는 두 개의 극한값을 이용해서 `$\displaystyle\lim _{\theta \rightarrow 0+} \dfrac{g(\theta)}{\theta \times f(\theta)}`$를 구합니다."

# Write cells in the exact order new shared strings were introduced in the
# source workbook revision.
$ws.Range("A209").Value = "y0019"
$ws.Range("B209").Value = $txt_B209
$ws.Range("B210").Value = $txt_B210
$ws.Range("C209").Value = "32111_x29"
$ws.Range("B211").Value = $txt_B211
$ws.Range("B212").Value = $txt_B212
$ws.Range("B213").Value = $txt_B213
$ws.Range("B214").Value = $txt_B214
$ws.Range("B215").Value = $txt_B215
$ws.Range("A210").Value = "y0020"
$ws.Range("A211").Value = "y0021"
$ws.Range("A212").Value = "y0022"
$ws.Range("A213").Value = "y0023"
$ws.Range("A214").Value = "y0024"
$ws.Range("A215").Value = "y0025"
$ws.Range("B216").Value = $txt_B216
$ws.Range("A216").Value = "y0026"

# Column C repeats the same problem id "32111_x29" for every new row.
$ws.Range("C210").Value = "32111_x29"
$ws.Range("C211").Value = "32111_x29"
$ws.Range("C212").Value = "32111_x29"
$ws.Range("C213").Value = "32111_x29"
$ws.Range("C214").Value = "32111_x29"
$ws.Range("C215").Value = "32111_x29"
$ws.Range("C216").Value = "32111_x29"

# Update the saved selection on the sheet view to A216.
$ws.Activate()
$ws.Range("A216").Select()
